$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update homework/lab scores for several students
$ws.Range("C8").Value = 2
$ws.Range("F11").Value = 2

$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 2

$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 2

$ws.Range("C30").Value = 2

# Update frozen pane top-left cell and active selection to reflect scroll position
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("F30").Select()
